# Adds a "UK" worksheet (Test Data for UK Market), cloned from the
# existing "Poland" sheet and placed right after it, matching the
# pattern used for every other per-country sheet in this workbook.

$wb = $excel.ActiveWorkbook
$poland = $wb.Worksheets.Item("Poland")

# Duplicate "Poland" and place the copy immediately after it - this is
# exactly what Excel's "Move or Copy... > Create a copy" does, and it
# carries over all formatting/styles/merged cells from the template sheet.
$poland.Copy($null, $poland)

# The newly inserted copy becomes the active sheet right after Copy().
$newSheet = $wb.ActiveSheet
$newSheet.Name = "UK"

# Fill in the two cells that differ from the Poland template.
# Order matters for shared-string table placement: the part number is
# written first so it lands right after "Poland Market", then the
# market label follows.
$newSheet.Range("B3").Value = "NGC-2741/T3357"
$newSheet.Range("B1").Value = "UK Market"

# Match the saved selection/active cell on the new sheet.
$newSheet.Range("B3").Select() | Out-Null
